# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" everywhere it is used
#    (Overview!E2:F3, and the "Status" column (C) on the zh-cn / de-de sheets).
# 2) Narrow the "zh-cn"/"de-de" status columns (Overview cols E/F, and column C
#    on the zh-cn/de-de detail sheets) from their old width down to match the
#    new narrower "Status" header width.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- 1. Update the status text wherever it appears -------------------------

$ws = $wb.Worksheets.Item("Overview")
foreach ($addr in @("E2", "F2", "E3", "F3")) {
    $cell = $ws.Range($addr)
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

foreach ($sheetName in @("zh-cn", "de-de")) {
    $detailWs = $wb.Worksheets.Item($sheetName)
    foreach ($addr in @("C2", "C3")) {
        $cell = $detailWs.Range($addr)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# --- 2. Narrow the now-shorter status columns -------------------------------

$ws.Columns.Item(5).ColumnWidth = 12.5   # Overview column E (zh-cn status)
$ws.Columns.Item(6).ColumnWidth = 12.5   # Overview column F (de-de status)

$wb.Worksheets.Item("zh-cn").Columns.Item(3).ColumnWidth = 12.5   # Status column
$wb.Worksheets.Item("de-de").Columns.Item(3).ColumnWidth = 12.5   # Status column
